$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "A1"  = 0.1621188509149718
    "A2"  = -0.0059999999365878409
    "A3"  = -0.0039999999447282164
    "A4"  = -0.0079999998966933106
    "A5"  = -0.0029999999422525292
    "A6"  = -0.001999999937282837
    "A7"  = -0.0099999998547999347
    "A8"  = -0.0099999998573969684
    "A9"  = -0.0019999999466344676
    "A10" = 0.049944486930726839
    "A11" = -0.002999999943033238
    "A12" = -0.0034999999374112356
    "A13" = -0.0034999999357481215
    "A14" = -0.007999999888904874
    "A15" = -0.00099999996177135841
    "A16" = -0.0019999999520283751
    "A17" = 0.011464490298361163
    "A18" = -0.0039999999328257374
    "A19" = -0.0039999999571187494
    "A20" = -0.0039999999536863839
    "A21" = -0.0039999999531827868
    "A22" = -0.0039999999527751129
    "A23" = -0.025799109927378439
    "A24" = 0.0085033845297486366
    "A25" = -0.019999999770527133
    "A26" = -0.0024999999292614206
    "A27" = -0.0024999999247956595
    "A28" = -0.0019999999072481955
    "A29" = -0.0069999998392304974
    "A30" = -0.0599999992821596
    "A31" = -0.0069999998266414565
    "A32" = -0.0099999997948643227
    "A33" = -0.0039999998568518436
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
